$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Element Types")
$ws.Range("B2:B5").Name = "Configuration"
$ws.Range("C2:C5").Name = "Configuration_Activity"
